$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NumberError")

# Mapping of row number -> new execution timestamp for column B (Date)
$rowTimestamps = @{
    2 = "Thu Nov 07 19:55:19 EST 2024"
    3 = "Thu Nov 07 19:55:33 EST 2024"
    4 = "Thu Nov 07 19:55:45 EST 2024"
    5 = "Thu Nov 07 19:55:56 EST 2024"
    6 = "Thu Nov 07 19:56:09 EST 2024"
    7 = "Thu Nov 07 19:56:23 EST 2024"
    8 = "Thu Nov 07 19:56:34 EST 2024"
    9 = "Thu Nov 07 19:56:47 EST 2024"
    10 = "Thu Nov 07 19:56:57 EST 2024"
    11 = "Thu Nov 07 19:57:09 EST 2024"
    12 = "Thu Nov 07 19:57:21 EST 2024"
    13 = "Thu Nov 07 19:57:33 EST 2024"
    14 = "Thu Nov 07 19:57:44 EST 2024"
    15 = "Thu Nov 07 19:57:55 EST 2024"
    16 = "Thu Nov 07 19:58:07 EST 2024"
    17 = "Thu Nov 07 19:58:17 EST 2024"
    18 = "Thu Nov 07 19:58:28 EST 2024"
    19 = "Thu Nov 07 19:58:39 EST 2024"
    20 = "Thu Nov 07 19:58:51 EST 2024"
    21 = "Thu Nov 07 19:59:01 EST 2024"
    22 = "Thu Nov 07 19:59:12 EST 2024"
    23 = "Thu Nov 07 19:59:23 EST 2024"
    24 = "Thu Nov 07 19:59:34 EST 2024"
    25 = "Thu Nov 07 19:59:45 EST 2024"
    26 = "Thu Nov 07 19:59:59 EST 2024"
    27 = "Thu Nov 07 20:00:11 EST 2024"
    28 = "Thu Nov 07 20:00:23 EST 2024"
    29 = "Thu Nov 07 20:00:34 EST 2024"
    30 = "Thu Nov 07 20:00:44 EST 2024"
    31 = "Thu Nov 07 20:00:55 EST 2024"
    32 = "Thu Nov 07 20:01:07 EST 2024"
    33 = "Thu Nov 07 20:01:18 EST 2024"
    34 = "Thu Nov 07 20:01:30 EST 2024"
    35 = "Thu Nov 07 20:01:42 EST 2024"
    36 = "Thu Nov 07 20:01:52 EST 2024"
    37 = "Thu Nov 07 20:02:04 EST 2024"
    38 = "Thu Nov 07 20:02:15 EST 2024"
    39 = "Thu Nov 07 20:02:27 EST 2024"
    40 = "Thu Nov 07 20:02:38 EST 2024"
    41 = "Thu Nov 07 16:26:49 EST 2024"
    42 = "Thu Nov 07 16:27:01 EST 2024"
    43 = "Thu Nov 07 17:45:07 EST 2024"
    44 = "Thu Nov 07 20:02:51 EST 2024"
    45 = "Thu Nov 07 20:03:01 EST 2024"
    46 = "Thu Nov 07 20:03:11 EST 2024"
    47 = "Thu Nov 07 20:03:22 EST 2024"
    48 = "Thu Nov 07 20:03:33 EST 2024"
    49 = "Thu Nov 07 20:03:43 EST 2024"
    50 = "Thu Nov 07 20:03:54 EST 2024"
    51 = "Thu Nov 07 20:04:06 EST 2024"
    52 = "Thu Nov 07 20:04:18 EST 2024"
    53 = "Thu Nov 07 20:04:29 EST 2024"
    54 = "Thu Nov 07 20:04:40 EST 2024"
    55 = "Thu Nov 07 20:04:51 EST 2024"
    56 = "Thu Nov 07 20:05:02 EST 2024"
    57 = "Thu Nov 07 20:05:13 EST 2024"
    58 = "Thu Nov 07 20:05:24 EST 2024"
    59 = "Thu Nov 07 20:05:35 EST 2024"
    60 = "Thu Nov 07 20:05:46 EST 2024"
    61 = "Thu Nov 07 20:05:56 EST 2024"
}

foreach ($row in $rowTimestamps.Keys) {
    $ws.Range("B$row").Value = $rowTimestamps[$row]
}

# Estate Tax rows (41-43) were removed from the RAD run: mark Result as DoNotRun
$estateRows = 41, 42, 43
foreach ($row in $estateRows) {
    $ws.Range("C$row").Value = "DoNotRun"
}

# Reflect the saved selection state: C41:C43 selected with C41 active
$ws.Activate()
$ws.Range("C41:C43").Select()
